# Add 2022-Q4 data:
#  - Insert a new "2022-Q4" worksheet (with fund holdings data) right before the
#    existing "2022-Q2" worksheet, so the tab order becomes: 总计, 2022-Q4, 2022-Q2.
#  - Update the "总计" (totals) summary sheet so that the 2022-Q4 quarter is the
#    new row 2, and the existing 2022-Q2 row is pushed down to row 3.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. "总计" sheet: move the existing 2022-Q2 summary row down to row 3, then
#    put the new 2022-Q4 summary in row 2.
# ---------------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.05

# Give the new A3 index cell the same look (bold/centered/bordered) as A2.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet, inserted before "2022-Q2" so it ends
#    up as the second tab.
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add($wsQ2)
$wsQ4.Name = "2022-Q4"

# Header row
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Row 2 - 中欧智能制造混合A
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "015143"
$wsQ4.Range("C2").Value = "中欧智能制造混合A"
$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "1.54"
$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "75.37"
$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "2.35"
$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.0362"
$wsQ4.Range("H2").Value = 10

# Row 3 - 中欧智能制造混合C
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "015144"
$wsQ4.Range("C3").Value = "中欧智能制造混合C"
$wsQ4.Range("D3").NumberFormat = "@"
$wsQ4.Range("D3").Value = "0.67"
$wsQ4.Range("E3").NumberFormat = "@"
$wsQ4.Range("E3").Value = "75.37"
$wsQ4.Range("F3").NumberFormat = "@"
$wsQ4.Range("F3").Value = "2.35"
$wsQ4.Range("G3").NumberFormat = "@"
$wsQ4.Range("G3").Value = "0.0157"
$wsQ4.Range("H3").Value = 10

# Match formatting (bold, centered, thin border) used by the other sheets'
# header rows / index column.
$wsTotal.Range("B1:D1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2").PasteSpecial(-4122)
$wsQ4.Range("A3").PasteSpecial(-4122)

$wsQ4.Range("A1").Select()
